$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "Muskan"
$ws.Range("B29").Value = "2021-01-15 06:23:48.161855"

$ws.Range("A30").Value = "Muskan Vaswan"
$ws.Range("B30").Value = "2021-01-15 07:07:40.004835"

$ws.Range("A31").Value = "Muskan Vaswan"
$ws.Range("B31").Value = "2021-01-15 07:08:33.727992"
